$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: add/modify explanatory notes ---
$ws.Range("B2").Value = "1 = white/caucasian, 2 = black, 3 = american indian, 4, asian/pacific islander, 5  = other (by modification) "
$ws.Range("C2").Value = "manually added"
$ws.Range("D2").Value = "0 = Male 1 = FEMALE,  9 = Missing or indeterminable (by modification) "

# --- Header row (row 1): STAMIN / STAMAX labels drop the word "requirement" ---
$ws.Range("G1").Value = "Total minimum statuary prison (STAMIN)"
$ws.Range("H1").Value = "Total maximum statuary prison (STAMAX)"

$ws.Range("J2").Value = "limit to under 1000 "
$ws.Rows.Item(2).RowHeight = 128

# --- Row 21: clear the long race-coding description, adjust row height ---
$ws.Range("B21").Value = ""
$ws.Rows.Item(21).RowHeight = 64

# --- Rows 22-30: clear the old race-code legend (1=White, 2=Black, ... NA=Missing) ---
$ws.Range("B22").Value = ""
$ws.Range("B23").Value = ""
$ws.Range("B24").Value = ""
$ws.Range("B25").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("B27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("B30").Value = ""

# --- View state: scroll position and selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("J11").Select()
